# Removing ocid entry from spreadsheet template. Fixes #46
#
# Every worksheet except "Activity" starts with a leading "ocid" column
# (column A). Remove that column from each of those sheets so the
# remaining columns shift left by one; the now-unused "ocid" shared
# string is pruned automatically on save.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "Classification",
    "Documents",
    "Event",
    "GrantProgramme",
    "Location",
    "Organization",
    "Transaction"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns("A").Delete()
}
